$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F values
$ws.Range("F1").Value = 0
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 17000

# Update selection to reflect new active cell (K8)
$ws.Range("K8").Select()
